# Adding missing files from previous commit
# Re-add the "Properties" worksheet (with its data) between "Leads" and "Opportunities".

$wb = $excel.ActiveWorkbook

$oppsSheet = $wb.Worksheets.Item("Opportunities")

# Insert a brand-new worksheet right before "Opportunities" (i.e. right after "Leads")
$props = $wb.Worksheets.Add($oppsSheet)
$props.Name = "Properties"

# Populate cells in the same order the values were originally entered so that
# the shared-strings table ends up in the same order as the source workbook.
$props.Range("C1").Value = "AskingPrice"
$props.Range("B1").Value = "State"
$props.Range("D1").Value = "Beds"
$props.Range("E1").Value = "Baths"
$props.Range("A2").Value = "Cambridge"
$props.Range("A1").Value = "City"
$props.Range("B2").Value = "MA"
$props.Range("F1").Value = "Broker"
$props.Range("F2").Value = "Victor Ochoa"
$props.Range("C2").Value = 450000
$props.Range("D2").Value = 3
$props.Range("E2").Value = 3

# Column sizing to match the authored layout.
$props.Columns.Item(1).ColumnWidth = 20.5625
$props.Columns.Item(2).ColumnWidth = 12.75
$props.Columns.Item(3).ColumnWidth = 10.1875
$props.Columns.Item(5).ColumnWidth = 7.3125
$props.Columns.Item(6).ColumnWidth = 12.4375

# The Baths column was formatted with 2 decimal places.
$props.Range("E2").NumberFormat = "0.00"

# Make the new Properties tab the active / selected sheet, zoomed to 130%,
# matching the tab-selection handoff from Leads seen in the source workbook.
$props.Activate()
$excel.ActiveWindow.Zoom = 130
